# Updates crypto price (D) and 1h volume change (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.889.38"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.361.39"
$ws.Range("E3").Value = "  -1.59%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'319.42"
$ws.Range("E5").Value = "  -4.70%  "
$ws.Range("D6").Value = "'107.34"
$ws.Range("E6").Value = "  +6.73%  "
$ws.Range("D7").Value = "'0.636"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.622"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").Value = "'41.31"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "'8.58"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "'1.01"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "'0.107"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "'15.94"
$ws.Range("E15").Value = "  -5.89%  "
$ws.Range("D16").Value = "2.714.40"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "2.352.40"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "42.846.47"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'0.0000107"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "'76.25"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'3.65"
$ws.Range("E22").Value = "  -6.65%  "
$ws.Range("D23").Value = "'266.62"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'2.32"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").Value = "'9.46"
$ws.Range("E25").Value = "  -7.70%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'11.44"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'23.44"
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "'36.76"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'167.78"
$ws.Range("E31").Value = "  -3.79%  "
$ws.Range("D32").Value = "'0.0909"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'6.00"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "'2.90"
$ws.Range("E34").Value = "  -6.96%  "
$ws.Range("D35").Value = "'0.131"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("D36").Value = "'0.119"
$ws.Range("E36").Value = "  +10.49%  "
$ws.Range("D37").Value = "'4.74"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'0.0362"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'3.86"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'2.73"
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("D41").Value = "'105.38"
$ws.Range("E41").Value = "  +12.20%  "
$ws.Range("D42").Value = "'1.52"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'0.240"
$ws.Range("E43").Value = "  +3.07%  "
$ws.Range("D44").Value = "'71.17"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("D47").Value = "'113.56"
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("D48").Value = "'5.54"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'9.16"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'76.25"
$ws.Range("E50").Value = "  +8.66%  "
$ws.Range("D51").Value = "'1.30"
$ws.Range("E51").Value = "  +0.87%  "
